# Bold + color (hex 2C3E50) highlighting for quantitative impact metrics
# (percentages, dollar amounts, counts) inside specific resume bullet
# paragraphs, matching the "hybrid bold + color highlighting" feature
# described in the commit message.

$d = $word.ActiveDocument
$metricColor = 5258796   # wdColor value for RGB(0x2C, 0x3E, 0x50) -> R + G*256 + B*65536

function Get-ParagraphByAnchor($doc, $mustContain, $mustNotContain) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $candidate = $doc.Paragraphs.Item($i)
        $text = $candidate.Range.Text
        if ($text.Contains($mustContain)) {
            if ([string]::IsNullOrEmpty($mustNotContain) -or -not $text.Contains($mustNotContain)) {
                return $candidate
            }
        }
    }
    return $null
}

function Highlight-Metric($paragraph, $searchText) {
    $range = $paragraph.Range
    $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Bold = $true
        $range.Font.Color = $metricColor
    }
    return $found
}

# 1) "Discovered systematic race coding errors ... demographic classification
#    accuracy from 23% to 64%" (Siege Analytics achievement bullet)
$para = Get-ParagraphByAnchor $d "Discovered systematic race coding errors" "KEY ACHIEVEMENTS"
Highlight-Metric $para "23%" | Out-Null
$para = Get-ParagraphByAnchor $d "Discovered systematic race coding errors" "KEY ACHIEVEMENTS"
Highlight-Metric $para "64%" | Out-Null

# 2) "Achieved 87% prediction accuracy for voter turnout vs. industry standard
#    of 71%, reducing polling error margins from ±4.2% to ±2.1%"
$para = Get-ParagraphByAnchor $d "reducing polling error margins" ""
Highlight-Metric $para "87%" | Out-Null
$para = Get-ParagraphByAnchor $d "reducing polling error margins" ""
Highlight-Metric $para "71%" | Out-Null
$para = Get-ParagraphByAnchor $d "reducing polling error margins" ""
Highlight-Metric $para "±4.2%" | Out-Null
$para = Get-ParagraphByAnchor $d "reducing polling error margins" ""
Highlight-Metric $para "±2.1%" | Out-Null

# 3) "Wrote RFP and analyzed bids from 1,200 vendors for research platform
#    development"
$para = Get-ParagraphByAnchor $d "Wrote RFP and analyzed bids from" ""
Highlight-Metric $para "1,200" | Out-Null

# 4) "Created comprehensive meta-analysis framework ... became the $400M
#    Polling Consortium Database at The Analyst Institute, now valued at $1B+"
$para = Get-ParagraphByAnchor $d "Polling Consortium Database" ""
Highlight-Metric $para "`$400M" | Out-Null
$para = Get-ParagraphByAnchor $d "Polling Consortium Database" ""
Highlight-Metric $para "`$1B" | Out-Null

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and
#    organizations $4.7M" (Key Achievements bullet)
$para = Get-ParagraphByAnchor $d "Algorithm reduced mapping costs by" ""
Highlight-Metric $para "73.5%" | Out-Null
$para = Get-ParagraphByAnchor $d "Algorithm reduced mapping costs by" ""
Highlight-Metric $para "`$4.7M" | Out-Null

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard
#    of 71%" (Key Achievements bullet, shorter variant without the ± margins)
$para = Get-ParagraphByAnchor $d "Achieved 87% prediction accuracy" "reducing polling"
Highlight-Metric $para "87%" | Out-Null
$para = Get-ParagraphByAnchor $d "Achieved 87% prediction accuracy" "reducing polling"
Highlight-Metric $para "71%" | Out-Null
